$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect before editing, then restore protection.
$ws.Unprotect()

# Update the confidential disclosure text's date (2021-05-06 -> 2021-05-07).
$ws.Range("A18").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-15.
$ws.Range("D2").Value2 = 0.05690919873797536
$ws.Range("E2").Value2 = 0.008482780876861495

$ws.Range("D3").Value2 = 0.02366384758734815
$ws.Range("E3").Value2 = 0.01164144353899887

$ws.Range("D4").Value2 = 0.03130451948005336
$ws.Range("E4").Value2 = 0.008862907788044438

$ws.Range("D5").Value2 = 0.0327507965444722
$ws.Range("E5").Value2 = 0.01636225266362246

$ws.Range("D6").Value2 = 0.03862248465362295
$ws.Range("E6").Value2 = 0.02873446097411847

$ws.Range("D7").Value2 = 0.01942152669216811
$ws.Range("E7").Value2 = 0.01276116897792456

$ws.Range("D8").Value2 = 0.004260122237369761
$ws.Range("E8").Value2 = 0.01801385681293333

$ws.Range("D9").Value2 = 0.006942129909210401
$ws.Range("E9").Value2 = 0.007558578987150399

$ws.Range("D10").Value2 = 0.07099547820983879
$ws.Range("E10").Value2 = 0.009423503325942661

$ws.Range("D11").Value2 = 0.07107418716572553
$ws.Range("E11").Value2 = 0.008859357696566983

$ws.Range("D12").Value2 = 0.1468394281023007
$ws.Range("E12").Value2 = -0.004931389365351646

$ws.Range("D13").Value2 = 0.3828502000525855
$ws.Range("E13").Value2 = 0.0003494975972040404

$ws.Range("D14").Value2 = 0.1143660806273293
$ws.Range("E14").Value2 = 0.002890521498253529

$ws.Range("E15").Value2 = 0.004097365234755346

# Restore sheet protection as it was before the edit.
$ws.Protect("D382")
